# Update latest output (run 134)
# Rewrites the pump-schedule optimisation results on the "Schedule" sheet
# and the corresponding detailed price/status series on the "Detailed" sheet.

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# ---- Schedule sheet: rows 2-5 (Start/Stop time, Duration, Volume, Cost, Unit Cost) ----
$wsSchedule.Range("A2").Value = 46042.02083333334
$wsSchedule.Range("B2").Value = 46042.1875
$wsSchedule.Range("C2").Value = 4
$wsSchedule.Range("D2").Value = 15.12
$wsSchedule.Range("E2").Value = 445.0682339999999
$wsSchedule.Range("F2").Value = 29.43572976190476

$wsSchedule.Range("A3").Value = 46042.29166666666
$wsSchedule.Range("C3").Value = 9
$wsSchedule.Range("D3").Value = 34.02
$wsSchedule.Range("E3").Value = -30.35730749999998
$wsSchedule.Range("F3").Value = -0.8923370811287473

$wsSchedule.Range("A4").Value = 46042.875
$wsSchedule.Range("C4").Value = 6
$wsSchedule.Range("D4").Value = 22.68
$wsSchedule.Range("E4").Value = 721.17638775
$wsSchedule.Range("F4").Value = 31.79790069444444

$wsSchedule.Range("E5").Value = -221.8333455
$wsSchedule.Range("F5").Value = -6.520674470899471

# ---- Detailed sheet: Pump_Status flips on 2023-12-19 early rows ----
$wsDetailed.Range("E2").Value = "OFF"
$wsDetailed.Range("E11").Value = "OFF"
$wsDetailed.Range("E12").Value = "OFF"
$wsDetailed.Range("E15").Value = "OFF"

# ---- Detailed sheet: rows 44-49 (historical/forecast boundary + status) ----
$wsDetailed.Range("E44").Value = "ON"

$wsDetailed.Range("B45").Value = 48.46067
$wsDetailed.Range("E45").Value = "ON"

$wsDetailed.Range("B46").Value = 56.98
$wsDetailed.Range("E46").Value = "ON"

$wsDetailed.Range("B47").Value = 57.3
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("E47").Value = "ON"

$wsDetailed.Range("B48").Value = 65.85254
$wsDetailed.Range("C48").Value = "historical"

$wsDetailed.Range("C49").Value = "historical"

# ---- Detailed sheet: rows 50-97 (updated price series) ----
$wsDetailed.Range("B50").Value = 64.8901
$wsDetailed.Range("B51").Value = 64.89
$wsDetailed.Range("B52").Value = 64.89
$wsDetailed.Range("B54").Value = 64.53478
$wsDetailed.Range("B55").Value = 64.89
$wsDetailed.Range("B57").Value = 65.64212999999999
$wsDetailed.Range("B58").Value = 68.87112999999999
$wsDetailed.Range("B59").Value = 73.20005
$wsDetailed.Range("B60").Value = 73.20005
$wsDetailed.Range("B61").Value = 79.95022
$wsDetailed.Range("B62").Value = 103.52886
$wsDetailed.Range("B63").Value = 77.94
$wsDetailed.Range("B64").Value = 35.88
$wsDetailed.Range("B65").Value = 0.66826
$wsDetailed.Range("B66").Value = -1.08193
$wsDetailed.Range("B67").Value = -5.97579
$wsDetailed.Range("B68").Value = -6.79084
$wsDetailed.Range("B69").Value = -7.41772
$wsDetailed.Range("B70").Value = -9.5
$wsDetailed.Range("B71").Value = -10
$wsDetailed.Range("B73").Value = -15.89865
$wsDetailed.Range("B74").Value = -22.10072
$wsDetailed.Range("B75").Value = -23.5
$wsDetailed.Range("B76").Value = -24.41017
$wsDetailed.Range("B78").Value = -27
$wsDetailed.Range("B79").Value = -25.94511
$wsDetailed.Range("B80").Value = -24.08764
$wsDetailed.Range("B81").Value = -22.86107
$wsDetailed.Range("B82").Value = -6.8
$wsDetailed.Range("B83").Value = -5.51
$wsDetailed.Range("B85").Value = 48.61802
$wsDetailed.Range("B86").Value = 55.37363
$wsDetailed.Range("B87").Value = 61.05003
$wsDetailed.Range("B89").Value = 79.95
$wsDetailed.Range("B90").Value = 73.37
$wsDetailed.Range("B91").Value = 65
$wsDetailed.Range("B92").Value = 57.31
$wsDetailed.Range("B93").Value = 59.66383
$wsDetailed.Range("B95").Value = 61.91795
$wsDetailed.Range("B96").Value = 63.33725
$wsDetailed.Range("B97").Value = 63.95727
